$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include #0")

# Update URL on Metadata sheet (row 2, B) to the new ValueSet URL
$wsMeta.Range("B2").Value = "http://hl7.org/fhir/ValueSet/immunization-status-reason"

# Update Date on Metadata sheet (row 8, B) to the new timestamp
$wsMeta.Range("B8").Value = "2025-06-25T06:29:04+01:00"

# Update System URI value on Include sheet (row 4, B) to match the same URL
# (this collapses the now-duplicate shared string)
$wsInclude.Range("B4").Value = "http://hl7.org/fhir/ValueSet/immunization-status-reason"
